# Updates the DTCP Results sheet (Ariyalur -> Chengalpattu dataset swap).
# Rows 2-11, columns A-J are replaced with the Chengalpattu "Building Plan"
# records; columns K (PDF link) and row 1 (headers) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "SWP/BPA/049180/2025"
$ws.Range("B2").Value = "Chengalpattu"
$ws.Range("C2").Value = "Building Plan"
$ws.Range("D2").Value = "'07/01/2025"
$ws.Range("E2").Value = "'13,779,223.00"
$ws.Range("F2").Value = "PLAN SHOWING THE PROPOSED CONSTRUCTION OF HIGH RISE COMMERCIAL BUILDING CONSIST OF BASEMENT + GROUND + 3 UPPER FLOORS COMPRISED IN SURVEY No: 107B, 108/1, 109/3, 109/5, 111/1A1 & 111/1B1 OF MELAKOTTAIYUR VILLAGE, VANDALUR TALUK, CHENGALPET DISTRICT."
$ws.Range("G2").Value = "For CASAGRAND MAGNUM PRIVATE LIMITED"
$ws.Range("H2").Value = "P. ELANCHEZIYAN, B.Arch, M.T.P Reg. Architect CA/2003/30729"

# Row 3
$ws.Range("A3").Value = "SWP/BPA/141104/2024"
$ws.Range("B3").Value = "Chengalpattu"
$ws.Range("C3").Value = "Building Plan"
$ws.Range("D3").Value = "'27/02/2025"
$ws.Range("E3").Value = "'13,478,056.00"
$ws.Range("F3").Value = "PLAN SHOWING THE PROPOSED CONSTRUCTION OF STILT + 5 FLOORS WITH 93 D.UNITS AFFORDABLE HOUSING RESIDENTIAL BUILDING AT S.NOS: 69/9D1A1 & 70/2A1, EAST COAST ROAD, KANATHUR REDDIKUPPAM VILLAGE, THIRUPORUR TALUK, CHENGALPET DISTRICT."
$ws.Range("G3").Value = "For SIDHARTH FOUNDATIONS AND HOUSING LIMITED"
$ws.Range("H3").Value = "G.A VAMSI VARMA K, B.Arch Regd..CA/2012/55059 4th Floor, Bootstart Co-working Rishab Arcade, Sanjay Nagar Main Road, Raj Mahal Villas 2nd Stage, Bengaluru, Karnataka-560001"
$ws.Range("J3").Value = "'9538974027"

# Row 4
$ws.Range("A4").Value = "SWP/BPA/208957/2024"
$ws.Range("B4").Value = "Chengalpattu"
$ws.Range("C4").Value = "Building Plan"
$ws.Range("D4").Value = "'27/02/2025"
$ws.Range("E4").Value = "'18,278,125.00"
$ws.Range("F4").Value = "PLAN SHOWING THE PROPOSED CONSTRUCTION OF INDUSTRIAL BUILDING OF PLOT NO: 6,7 & 8 IN (SURVEY NUMBER 10/1pt, 2pt, 4pt, 5pt, 6pt, 11/1Apt, 1Bpt, 1C, 1D, 1E, 14pt, 15/1pt, 2Apt, 2Bpt, 2Cpt, 16pt, 17/1Fpt, 1Gpt, 1H, 1I, 1J, 2A2Bpt, 2A3Apt, 2A3Bpt, 2A4pt, 2B1, 2B2, 2B3, 2B4, 2B5, 3A, 3B, 4A, 4B, 4C, 4D AND 4Ept.) PANCHANTHIRUTHI VILLAGE, THIRUPORUR TALUK, CHENGALPATTU DISTRICT."
$ws.Range("G4").Value = "For HEADWAY PREMIER INDUSTRIES PRIVATE LIMITED Authorised Signatory"
$ws.Range("H4").Value = "A.N. RAVICHANDRAN, DCE, B.E. (Civil), ME, AIV Reg. Registered Engineer Grade-1B Reg: ENG/2022/70/10000 Structural Engineer Grade-1B Reg: ENG/2022/70/10005"
$ws.Range("I4").Value = "mirrahbuilders@gmail.com"
$ws.Range("J4").Value = "+91 9884147407"

# Row 5
$ws.Range("A5").Value = "SWP/BPA/223491/2024"
$ws.Range("B5").Value = "Chengalpattu"
$ws.Range("C5").Value = "Building Plan"
$ws.Range("D5").Value = "'28/01/2025"
$ws.Range("E5").Value = "'3,133,961.00"
$ws.Range("F5").Value = "PLAN SHOWING THE PROPOSED CONSTRUCTION OF DATA COMMUNICATION CENTER (COMMERCIAL BUILDING) AT PLOT NO: 42 IN APPROVED LAY OUT NO: L.P/D.T.C.P NO:76/2021 CHENGALPATTU. COMPRISED IN SURVEY NO.36/3B(Pt.), 4A(Pt.), 4B(Pt.), 5A(Pt.), 5B1(Pt.), 5B2(Pt.), 8B(Pt.), 9A(Pt.), 37/1A(Pt.), 1B(Pt.), 2A(Pt.), 2B(Pt.), 2B2(Pt.), 2B3(Pt.), 2B4(Pt.), 2B5(Pt.), 2C(Pt.), 3(Pt.), 4A1(Pt.), 4A2(Pt.), 4C, 4D1(Pt.), 4D2(Pt.), 4E1(Pt.), 5B(Pt.), 5C, 5D(Pt.), 40/1A1A(Pt.), 1A2, 1C(Pt.), 1D1, 1D2(Pt.), 1D3(Pt.), 2B1(Pt.), 2B2, 2C(Pt.), 3B, 3C1(Pt.), 3C2, 4A, 4B, 4C, 5A(Pt.), 5B(Pt.), 41/1A1(Pt.), 1B, 1C, 1D(Pt.), 2(Pt.), 3, 4, 5, 42/1(Pt.), 2A, 2B1A(Pt.), 2B1B(Pt.), 52/1A1A6(Pt.), 1B(Pt.) OF KUNNAPATTU VILLAGE, THIRUPORUR TALUK, KUNNAPATTU PANCHAYAT THIRUPORUR PANCHAYAT UNION CHENGALPATTU DISTRICT."
$ws.Range("G5").Value = "For TATA COMMUNICATIONS LIMITED"
$ws.Range("H5").Value = "K. JAMAL MOHIDEEN, B.Arch Registered Architect-CA/1993/7686 New .17, Old .5, 4th Street, Ashok Nagar, Chennai-600 083"
$ws.Range("I5").Value = "cparchitect@gmail.com"
$ws.Range("J5").Value = "98400 69523"

# Row 6
$ws.Range("A6").Value = "SWP/BPA/233615/2024"
$ws.Range("B6").Value = "Chengalpattu"
$ws.Range("C6").Value = "Building Plan"
$ws.Range("D6").Value = "'14/02/2025"
$ws.Range("E6").Value = "'3,971,039.00"
$ws.Range("F6").Value = "PLAN SHOWING THE EXISTING CONSTRUCTED BUILDING & PROPOSED CONSTRUCTION OF INDUSTRIAL BUILDING IN S.NO:152 / 1(PART), AS PER PATTA S.NO:152/IC, PLOT NO:AIFS, MMDA INDUSTRIAL COMPLEX AT MARAMALAINAGAR, GUDDALUR VILLAGE , CHENGALPET TALUK & DISTRICT."
$ws.Range("G6").Value = "Ar K Rajeshkumar"
$ws.Range("H6").Value = "Dhinesh. T Architect CA/2019/49592 Reg TTLP/KR/AP/2019/05/001 .1/4, 1st Street, Kasturbai Nagar, Adyar, Chennai 600 020"
$ws.Range("I6").Value = "dhinesh.architect@gmail.com"
$ws.Range("J6").Value = "+91 7378038000 | +91 99600 66666"

# Row 7
$ws.Range("A7").Value = "SWP/BPA/240909/2024"
$ws.Range("B7").Value = "Chengalpattu"
$ws.Range("C7").Value = "Building Plan"
$ws.Range("D7").Value = "'20/01/2025"
$ws.Range("E7").Value = "'3,696,915.00"
$ws.Range("F7").Value = "PLAN SHOWING THE REVISED AND EXPANSION OF M/s. ARMSTRONG INTERNATIONAL PVT.LTD. BUILDING AT PLOT NO - P-46 AT MAHINDRA WORLD CITY, CHENGALPATTU, COMPRISED IN S.NO:165/3, 171/1, 171/2, 172/1, 191/3, 193/1,193/2, 197/1, 198/1, 198/2, 199/1, 199/2. OF ANJUR VILLAGE, MAHINDRA WORLD CITY, CHENGALPATTU TALUK, CHENGALPATTU DISTRICT."
$ws.Range("G7").Value = "For ARMSTRONG INTERNATIONAL Pvt LTD, Director - Finance."
$ws.Range("H7").Value = "Dhinesh.T, Architect CA/2009/49367 Reg. TVLR/RCA/Reg2/05/001 10/51, 1st Street, Kasturibai Nagar, Alwarthirunagar, Chennai 600 054"
$ws.Range("I7").Value = "dhinesh.architect@gmail.com"
$ws.Range("J7").Value = "+91 73730 83006 | +91 69003 66666"

# Row 8
$ws.Range("A8").Value = "SWP/BPA/270950/2025"
$ws.Range("B8").Value = "Chengalpattu"
$ws.Range("C8").Value = "Building Plan"
$ws.Range("D8").Value = "'03/01/2025"
$ws.Range("E8").Value = "'72,764,713.00"
$ws.Range("F8").Value = "PLAN SHOWING THE PROPOSED CONSTRUCTION OF GROUP DEVELOPMENT FOR AFFORDABLE RESIDENTIAL APARTMENTS HIGH RISE BUILDING IN S.F.NO.18/1, 2, 3, 8A, 8B, 8C, 10; 19/1, 2, 3, 4, 6B; 23/1A, 24/1, 25/1A, 1B, 1C, 12C3B; 30/5C1A, 6B, 7A, 8, 9, 10, 12C1A, 12C1B, 12C2, 12C3, 13A1, 13B; 31/1C1A, 11,13 AT KANATHUR REDDYKUPPAM VILLAGE, TIRUPORUR TALUK, CHENGALPATTU DISTRICT."
$ws.Range("G8").Value = "For Voora Property Developers (P) Ltd. Managing Director"
$ws.Range("H8").Value = "M/s. VERNEKAR ASSOCIATES PVT. LTD. . 8, 1st 'A' Main, Stag Extension ST Bed Koramangala 4th Block BANGALORE-560095 .: 40609800, Fax: 40609810"

# Row 9
$ws.Range("A9").Value = "SWP/BPA/325413/2024"
$ws.Range("B9").Value = "Chengalpattu"
$ws.Range("C9").Value = "Building Plan"
$ws.Range("D9").Value = "'21/01/2025"
$ws.Range("E9").Value = "'3,323,788.00"
$ws.Range("F9").Value = "PLAN SHOWING THE PROPOSED CONSTRUCTION OF INDUSTRIAL BUILDING (FORMULATIONS OF PHARMACEUTICALS) AT SIDCO INDUSTRIAL ESTATE, PLOT NO.53,54,55,82,83 & 84, COMPRISED IN OLD SURVEY NO.237/10A1A1A, NEW SURVEY NO.237/27 IN ALATHUR VILLAGE, THIRUPORUR TALUK, CHENGALPATTU DISTRICT."
$ws.Range("G9").Value = "For ORCHID PHARMA LIMITED"
$ws.Range("H9").Value = "B.MALLIKA, ME., (PhD) Registered Engineer Grade-II RE-II / KKDTCP / 017 / 2022 84, Kaalimann Koli Street, Pasumpon Nagar, Madurai-625003"

# Row 10
$ws.Range("A10").Value = "SWP/BPA/358599/2024"
$ws.Range("B10").Value = "Chengalpattu"
$ws.Range("C10").Value = "Building Plan"
$ws.Range("D10").Value = "'28/01/2025"
$ws.Range("E10").Value = "'249,081.00"
$ws.Range("F10").Value = "PLAN SHOWING OFFICE BUILDING`" IN STILT FLOOR & HEAD ROOM (279 SQ.M B) GUDUVANCHERI VILLAGE, CHENGALPATTU DISTRICT."
$ws.Range("G10").Value = "APPLICANT"
$ws.Range("H10").Value = "Er. U. UMESHRAJ, B.E., Registered Engineer Grade-II (RE) RE/CA/AUA/2024/031 .32, Nasri 1st Cross Street, Ezhil Nagar Ext., Thorapadi Post, Vellore - 632 002"
$ws.Range("I10").Value = "umeshraj1487@gmail.com"
$ws.Range("J10").Value = "99440 74734"

# Row 11
$ws.Range("A11").Value = "SWP/BPA/360192/2025"
$ws.Range("B11").Value = "Chengalpattu"
$ws.Range("C11").Value = "Building Plan"
$ws.Range("D11").Value = "'30/01/2025"
$ws.Range("E11").Value = "'49,619,230.00"
$ws.Range("F11").Value = "PLAN SHOWING THE PROPOSED ADDITIONAL CONSTRUCTION OF RESIDENTIAL HIGHRISE GROUP DEVELOPMENT BUILDING CONSISTS OF 3 BLOCKS, BLOCK 3 & 6 WITH COMMON CLUB HOUSE AT BLOCKS 1, 2, & 5. BLOCK 3: GROUND+13 UPPER FLOORS (25 DU); BLOCK 6: GROUND+18 UPPER FLOORS (252 DU); BLOCK 6 CLUB HOUSE: GROUND + 3 FLOORS WITH SWIMMING POOL ON THE GROUND FLOOR, TOTAL 500 DU COMPOSED IN S.NOS. 110/2A1B, 110/2A1C, 110/2A1D, 110/2A1E, 111/1C, 111/1D, 111/1E, 111/1F, 111/1G, 111/2A1, 111/2C, 112/2A2, 112/2A1, 114/C, 114/2, 117/2B, 117/2C, 117/4, 117/5, 118/1D, 118/2, 119/1, 119/2B, 119/2A, 119/4, 119/5, 122/2A, 122/1, 122/3, 122/4, 122/5, 122/6A, 122/6B, 124/1, 124/2, 124/3, 124/4, 124/5A, 124/5B, 124/5C, 124/5D, MELAKOTTAIYUR VILLAGE, VANDALUR TALUK, CHENGALPATTU DISTRICT."
$ws.Range("G11").Value = "For CASAGRAND ANCHOR PRIVATE LIMITED Authorised Signatory"
$ws.Range("H11").Value = "JEYAVEL GANDHAM M.S.S B.Arch. Registered Architect - CA/2008/42210"
